$d = $word.ActiveDocument

# Locate the "Nato/a: ..." paragraph inside the first table (Dati personali dello stagista).
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(2, 1)
$cellRange = $cell.Range
$cellStart = $cellRange.Start

# Replace the long filler of ellipses after "residente in: " with the new
# {S_RESIDENZA} placeholder while the paragraph still holds a single run
# (this keeps Find/Replace from merging/leaving stray runs).
$tailRange = $d.Range($cellStart + 35, $cellStart + 67)
$tailRange.Find.Execute("……………………………………………", $false, $false, $false, $false, $false, $true, 1, $false, "…{S_RESIDENZA}…", 2) | Out-Null

# Split the paragraph's single run into five runs matching the authored
# template: "Nato/a: …{S_" | "NATOA}…" | " il …{S_" | "NATOIL}…" |
# " residente in: …{S_RESIDENZA}…" -- toggling a character property on and
# back off forces Word to break the run at these boundaries without
# altering the final run formatting.
$split1 = $d.Range($cellStart + 12, $cellStart + 19)
$split1.Bold = 1
$split1.Bold = 0

$split2 = $d.Range($cellStart + 19, $cellStart + 27)
$split2.Bold = 1
$split2.Bold = 0

$split3 = $d.Range($cellStart + 27, $cellStart + 35)
$split3.Bold = 1
$split3.Bold = 0

Write-Output "Residenza field added to Nato/a paragraph"
